$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (Excel serial 45185 = 2023-09-16)
# that needs to be updated to serial 45204 (2023-10-05) for every data row
# (rows 2 through 173). Use the raw Excel serial number so no time-of-day
# fraction is introduced.
$lastRow = 173
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 45204
    }
}
